# Auto-update draw results: append the 2025-12-02 "Pick 4" draw as a new
# row (77) at the bottom of the results table on the active sheet, mirroring
# the existing rows' layout (Date, Game, Phase, Result, InsertedAt).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 77

# Column A ("Date") and column C ("Phase") hold digit-only strings that Excel's
# smart-entry would otherwise reinterpret as a date serial / plain number, so
# mark them as Text ("@") before writing, then put the style back to "Normal"
# so no stray per-cell formatting is left behind (matches the rest of the
# sheet, which carries no explicit cell styles).
$dateCell = $ws.Cells.Item($row, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025-12-02"
$dateCell.Style = "Normal"

$ws.Cells.Item($row, 2).Value = "Pick 4"

$phaseCell = $ws.Cells.Item($row, 3)
$phaseCell.NumberFormat = "@"
$phaseCell.Value = "251202"
$phaseCell.Style = "Normal"

$ws.Cells.Item($row, 4).Value = "1-1-0-1"
$ws.Cells.Item($row, 5).Value = "2025-12-02T21:45:09.960+04:00"
